# Adds the 2024/11/22 data column (BW) to the "合成確率" sheet, one column to
# the right of the current last column (BV, 2024/11/21), mirroring the layout
# of all the other daily columns: a text date header in row 1 plus 52 numeric
# values in rows 2-53, each carrying the same conditional-style look (plain /
# yellow-highlight / blue-highlight) as its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteSpecial constants used below.
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# --- Header cell BW1 ------------------------------------------------------
# A plain  Range("BW1").Value = "2024/11/22"  gets auto-recognised as a date
# literal and rewritten into a date-serial cell with a brand new number
# format, so instead we compute the literal text in a scratch cell via a
# formula that evaluates to a string, copy only the VALUE of that over to
# BW1 (landing it as plain text, same as the other date headers), clear the
# scratch cell again, and finally copy the FORMAT of the existing header
# cell (A1) onto BW1 so it picks up the same style as every other header.
$ws.Range("ZZ1").Formula = "=""2024/11/22"""
$ws.Range("ZZ1").Copy()
$ws.Range("BW1").PasteSpecial($xlPasteValues)
$ws.Range("ZZ1").Value = ""

$ws.Range("A1").Copy()
$ws.Range("BW1").PasteSpecial($xlPasteFormats)

# --- Data rows 2-53 --------------------------------------------------------
# Style 1 = plain (メイリオ, no fill), 2 = yellow highlight, 3 = blue highlight.
# Donor cells already carrying each style are reused via copy/paste-special
# so no new style/fill entries get created in the workbook.
$donors = @{ 1 = "A2"; 2 = "D2"; 3 = "N2" }

$data = @(
    @{ Row = 2; Style = 3; Value = 126.5 },
    @{ Row = 3; Style = 3; Value = 132.1 },
    @{ Row = 4; Style = 1; Value = 163.5 },
    @{ Row = 5; Style = 1; Value = 151.2 },
    @{ Row = 6; Style = 1; Value = 187 },
    @{ Row = 7; Style = 1; Value = 323 },
    @{ Row = 8; Style = 1; Value = 167.7 },
    @{ Row = 9; Style = 1; Value = 242.1 },
    @{ Row = 10; Style = 1; Value = 151.6 },
    @{ Row = 11; Style = 3; Value = 131.2 },
    @{ Row = 12; Style = 1; Value = 161.7 },
    @{ Row = 13; Style = 1; Value = 172.7 },
    @{ Row = 14; Style = 3; Value = 138.8 },
    @{ Row = 15; Style = 1; Value = 141.9 },
    @{ Row = 16; Style = 1; Value = 190 },
    @{ Row = 17; Style = 2; Value = 122 },
    @{ Row = 18; Style = 1; Value = 144.2 },
    @{ Row = 19; Style = 1; Value = 152.5 },
    @{ Row = 20; Style = 1; Value = 144.1 },
    @{ Row = 21; Style = 1; Value = 154.1 },
    @{ Row = 22; Style = 3; Value = 131.7 },
    @{ Row = 23; Style = 1; Value = 165.9 },
    @{ Row = 24; Style = 1; Value = 241.7 },
    @{ Row = 25; Style = 1; Value = 153.7 },
    @{ Row = 26; Style = 1; Value = 243.1 },
    @{ Row = 27; Style = 1; Value = 271.1 },
    @{ Row = 28; Style = 1; Value = 143.1 },
    @{ Row = 29; Style = 1; Value = 250.4 },
    @{ Row = 30; Style = 1; Value = 159.7 },
    @{ Row = 31; Style = 1; Value = 151.9 },
    @{ Row = 32; Style = 1; Value = 141.7 },
    @{ Row = 33; Style = 2; Value = 112.6 },
    @{ Row = 34; Style = 1; Value = 177.5 },
    @{ Row = 35; Style = 1; Value = 148.6 },
    @{ Row = 36; Style = 2; Value = 114.6 },
    @{ Row = 37; Style = 1; Value = 332.6 },
    @{ Row = 38; Style = 2; Value = 121.8 },
    @{ Row = 39; Style = 3; Value = 125.6 },
    @{ Row = 40; Style = 1; Value = 156.1 },
    @{ Row = 41; Style = 1; Value = 173.5 },
    @{ Row = 42; Style = 1; Value = 171.9 },
    @{ Row = 43; Style = 1; Value = 140.5 },
    @{ Row = 44; Style = 1; Value = 219.5 },
    @{ Row = 45; Style = 1; Value = 142.5 },
    @{ Row = 46; Style = 1; Value = 158.8 },
    @{ Row = 47; Style = 1; Value = 249.8 },
    @{ Row = 48; Style = 1; Value = 148.3 },
    @{ Row = 49; Style = 1; Value = 221 },
    @{ Row = 50; Style = 2; Value = 117.2 },
    @{ Row = 51; Style = 3; Value = 129.6 },
    @{ Row = 52; Style = 3; Value = 129.2 },
    @{ Row = 53; Style = 1; Value = 178.4 }
)

foreach ($item in $data) {
    $donor = $donors[$item.Style]
    $target = "BW" + $item.Row
    $ws.Range($donor).Copy()
    $ws.Range($target).PasteSpecial($xlPasteFormats)
    $ws.Range($target).Value = $item.Value
}

# --- Column width -----------------------------------------------------------
# 11.2 characters of COM ColumnWidth round-trips to the same OOXML
# width="12" used by every other data column on this sheet.
$ws.Columns.Item(75).ColumnWidth = 11.2
